$wb = $excel.ActiveWorkbook

$wsGrilla = $wb.Worksheets.Item("grilla de pruebas")
$wsCalc = $wb.Worksheets.Item("CALCULADORA")

# --- "grilla de pruebas" sheet edits ---
$wsGrilla.Range("B1").Value = "BUY"
$wsGrilla.Range("B3").Value = "534.95"
$wsGrilla.Range("F3").Value = "0.06475"

# B9 is formatted as Text ("@") but must hold a real numeric value.
$wsGrilla.Range("B9").NumberFormat = "General"
$wsGrilla.Range("B9").Value = 10
$wsGrilla.Range("B9").NumberFormat = "@"

$wsGrilla.Range("C9").Value = "VA 6. Probando 10 automatico y X manual"

# --- "CALCULADORA" sheet edits ---
$wsCalc.Range("A2").Value = "BUY"

# B3, E4, F4 are formatted as Text ("@") but must hold real numeric
# values (as in the source workbook). Temporarily switch to General
# so the typed value is stored as a number, then restore the Text format.
$wsCalc.Range("B3").NumberFormat = "General"
$wsCalc.Range("B3").Value = 538.63
$wsCalc.Range("B3").NumberFormat = "@"

$wsCalc.Range("E3").Value = 1678
$wsCalc.Range("F3").Value = 0.006075

$wsCalc.Range("E4").NumberFormat = "General"
$wsCalc.Range("E4").Value = 0
$wsCalc.Range("E4").NumberFormat = "@"

$wsCalc.Range("F4").NumberFormat = "General"
$wsCalc.Range("F4").Value = 0.006031
$wsCalc.Range("F4").NumberFormat = "@"

$wsCalc.Range("J6").Value = 0.006115
$wsCalc.Range("K6").Formula = "=H6/100*SUM(G3:G4)*J6"

# --- selections / active sheet ---
$wsCalc.Activate()
$wsCalc.Range("D21").Select()

$wsGrilla.Activate()
$wsGrilla.Range("B10").Select()
